# Add the "second half" of the Bejarano (Palau) bite-rate observations.
# Rows 26-40 already existed as placeholder rows (Paper/Location/Unit/Person
# filled in, but Notes/Bite rate/Species empty) - they need their values
# filled in, plus 4 brand-new rows (41-44) for the rest of the data set.
# The existing Bellwood_Choat block (old rows 41-64) shifts down by 4 rows
# to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the old rows 41:64 (Bellwood_Choat data) down by 4 rows so rows
# 41-44 are free for the new Bejarano observations.
$ws.Rows("41:44").Insert()

# NOTE: this COM-interop PowerShell engine only binds positional
# arguments reliably (named parameters silently fail to bind), so
# Set-RowData is called positionally below.
function Set-RowData {
    param($Row, $Notes, $BiteRate, $Species)
    $ws.Range("A$Row").Value = "Bejarano"
    $ws.Range("B$Row").Value = "Palau"
    $ws.Range("C$Row").Value = $Notes
    $ws.Range("D$Row").Value = "bites per hour per m^2"
    $ws.Range("E$Row").Value = $BiteRate
    $ws.Range("E$Row").HorizontalAlignment = -4108
    $ws.Range("F$Row").Value = $Species
    $ws.Range("G$Row").Value = "Jan"
}

Set-RowData 26 "low wave"    1.1235955056179774   "Scarus oviceps"
Set-RowData 27 "medium wave" 16.256759113901971   "Scarus psittacus"
Set-RowData 28 "high wave"   77.205040091638026   "Scarus psittacus"
Set-RowData 29 "medium wave" 0.44444444444444448  "Siganus puellus"
Set-RowData 30 "medium wave" 3.4090909090909092   "Siganus punctatus"
Set-RowData 31 "high wave"   2.9166666666666665   "Siganus punctatus"
Set-RowData 32 "low wave"    28.886786509972211   "Scarus schlegeli"
Set-RowData 33 "medium wave" 18.909698840254396   "Scarus schlegeli"
Set-RowData 34 "high wave"   152.53045605446582   "Scarus schlegeli"
Set-RowData 35 "low wave"    14.166666666666668   "Scarus spinus"
Set-RowData 36 "medium wave" 103.98518445839875   "Scarus spinus"
Set-RowData 37 "high wave"   221.11810949087476   "Scarus spinus"
Set-RowData 38 "low wave"    5.0493986254295535   "Siganus vulpinus"
Set-RowData 39 "medium wave" 18.315217391304348   "Siganus vulpinus"
Set-RowData 40 "high wave"   0.33333333333333331  "Siganus vulpinus"
Set-RowData 41 "low wave"    328.12092005062789   "Zebrasoma scopas"
Set-RowData 42 "medium wave" 253.73992677916567   "Zebrasoma scopas"
Set-RowData 43 "high wave"   6.4016064257028118   "Zebrasoma scopas"
Set-RowData 44 "low wave"    19.2053264604811     "Zebrasoma veliferum"

# Update the view to match where the author ended up working.
$ws.Range("E44").Select()
$win = $excel.ActiveWindow
$win.Zoom = 80
